$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Insert()

$ws.Range("E7").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats

$d7 = $ws.Range("D7")
$d7.Value = 43465

Write-Host "done"
